# Apply the edit described by the diff:
#  - On worksheet "Method1", cells B38:B77 change their looked-up standard
#    name from "dCL 61:1" to "dCL 80:4". The C/D/E columns are VLOOKUP-driven
#    formulas referencing StdInfo, so they recalculate automatically.
#  - The active selection on "Method1" moves from B4 to B6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Method1")

$ws.Range("B38:B77").Value = "dCL 80:4"

[void]$ws.Range("B6").Select()
